$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet to match the file name
$ws.Name = "sandwiches_greenmountain"

# Extend the table to cover the new column (G) and drop the empty trailing row (5)
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:G4"))

# Header row - name the new column
$ws.Range("G1").Value = "LeaveEmpty"

# Row 2 - Club Sandwich
$ws.Range("A2").Value = "Club Sandwich"
$ws.Range("D2").Value = "Island City Bread"
$ws.Range("E2").Value = "GF, BC"
$ws.Range("F2").Value = "cobbsalad"

# Row 3 - Turkey Swiss Wrap
$ws.Range("A3").Value = "Turkey Swiss Wrap"
$ws.Range("E3").Value = "GF,VEG, BC"
$ws.Range("F3").Value = "kalecaesar"

# Row 4 - Falafel Wrap
$ws.Range("A4").Value = "Falafel Wrap"
$ws.Range("C4").Value = "No known priority allergens"
$ws.Range("D4").Value = "Zorba's Tzatziki BCfresh Tomatoes"
$ws.Range("E4").Value = "VGN, BC, GF, DF"
$ws.Range("F4").Value = "housesalad"

# Match final selection state
[void]$ws.Range("D4").Select()
